# Scheduled-runner refresh of market-price-derived columns
# (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ])
# across the per-job "Alexander_Profits" worksheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR). Only numeric H:N cells move; A:G (leve metadata)
# are left untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 310
$ws.Range("I33").Value = 187.5
$ws.Range("J33").Value = 800
$ws.Range("K33").Value = 187.5
$ws.Range("L33").Value = 800
$ws.Range("M33").Value = 41.5
$ws.Range("N33").Value = -1258
$ws.Range("H40").Value = 1354.8937
$ws.Range("I40").Value = 1189.2307
$ws.Range("J40").Value = 1560
$ws.Range("K40").Value = 1189.2307
$ws.Range("L40").Value = 1560
$ws.Range("M40").Value = -1014.2307
$ws.Range("N40").Value = -1910
$ws.Range("H45").Value = 4271.4287
$ws.Range("J45").Value = 4271.4287
$ws.Range("L45").Value = 12814.2861
$ws.Range("N45").Value = -13198.2861
$ws.Range("H116").Value = 5114.7617
$ws.Range("I116").Value = 3786.25
$ws.Range("J116").Value = 6886.1113
$ws.Range("K116").Value = 3786.25
$ws.Range("L116").Value = 6886.1113
$ws.Range("M116").Value = -344.25
$ws.Range("N116").Value = -13770.1113
$ws.Range("H137").Value = 1638.7046
$ws.Range("I137").Value = 1250.1428
$ws.Range("J137").Value = 2318.6875
$ws.Range("K137").Value = 3750.4284
$ws.Range("L137").Value = 6956.0625
$ws.Range("M137").Value = -1200.4284
$ws.Range("N137").Value = -12056.0625
$ws.Range("H138").Value = 2439.0476
$ws.Range("I138").Value = 2751.2727
$ws.Range("J138").Value = 2328.258
$ws.Range("K138").Value = 8253.8181
$ws.Range("L138").Value = 6984.773999999999
$ws.Range("M138").Value = -3113.8181
$ws.Range("N138").Value = -17264.774

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2115.9167
$ws.Range("I2").Value = 1364.2778
$ws.Range("J2").Value = 4370.8335
$ws.Range("K2").Value = 1364.2778
$ws.Range("L2").Value = 4370.8335
$ws.Range("M2").Value = -1251.2778
$ws.Range("N2").Value = -4596.8335
$ws.Range("H32").Value = 26195.803
$ws.Range("I32").Value = 5854.518
$ws.Range("J32").Value = 140107
$ws.Range("K32").Value = 5854.518
$ws.Range("L32").Value = 140107
$ws.Range("M32").Value = -5567.518
$ws.Range("N32").Value = -140681
$ws.Range("H61").Value = 1695.75
$ws.Range("I61").Value = 1501.9445
$ws.Range("K61").Value = 1501.9445
$ws.Range("M61").Value = -1289.9445
$ws.Range("H63").Value = 2659.4119
$ws.Range("I63").Value = 1523.3334
$ws.Range("J63").Value = 3937.5
$ws.Range("K63").Value = 1523.3334
$ws.Range("L63").Value = 3937.5
$ws.Range("M63").Value = -837.3334
$ws.Range("N63").Value = -5309.5
$ws.Range("H66").Value = 2659.4119
$ws.Range("I66").Value = 1523.3334
$ws.Range("J66").Value = 3937.5
$ws.Range("K66").Value = 7616.666999999999
$ws.Range("L66").Value = 19687.5
$ws.Range("M66").Value = -4184.666999999999
$ws.Range("N66").Value = -26551.5
$ws.Range("H102").Value = 50001292
$ws.Range("I102").Value = 1401.3334
$ws.Range("J102").Value = 200000960
$ws.Range("K102").Value = 1401.3334
$ws.Range("L102").Value = 200000960
$ws.Range("M102").Value = 220.6666
$ws.Range("N102").Value = -200004204
$ws.Range("H116").Value = 2115.9167
$ws.Range("I116").Value = 1364.2778
$ws.Range("J116").Value = 4370.8335
$ws.Range("K116").Value = 1364.2778
$ws.Range("L116").Value = 4370.8335
$ws.Range("M116").Value = 929.7221999999999
$ws.Range("N116").Value = -8958.833500000001
$ws.Range("H136").Value = 1695.75
$ws.Range("I136").Value = 1501.9445
$ws.Range("K136").Value = 4505.833500000001
$ws.Range("M136").Value = -1955.833500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2115.9167
$ws.Range("I3").Value = 1364.2778
$ws.Range("J3").Value = 4370.8335
$ws.Range("K3").Value = 1364.2778
$ws.Range("L3").Value = 4370.8335
$ws.Range("M3").Value = -1250.2778
$ws.Range("N3").Value = -4598.8335
$ws.Range("H58").Value = 23692
$ws.Range("J58").Value = 23692
$ws.Range("L58").Value = 23692
$ws.Range("N58").Value = -24280
$ws.Range("H80").Value = 504.625
$ws.Range("I80").Value = 298.66666
$ws.Range("J80").Value = 552.1539
$ws.Range("K80").Value = 298.66666
$ws.Range("L80").Value = 552.1539
$ws.Range("M80").Value = 699.33334
$ws.Range("N80").Value = -2548.1539
$ws.Range("H83").Value = 504.625
$ws.Range("I83").Value = 298.66666
$ws.Range("J83").Value = 552.1539
$ws.Range("K83").Value = 1493.3333
$ws.Range("L83").Value = 2760.7695
$ws.Range("M83").Value = 3498.6667
$ws.Range("N83").Value = -12744.7695
$ws.Range("H99").Value = 2038.6842
$ws.Range("I99").Value = 1102.5
$ws.Range("J99").Value = 4660
$ws.Range("K99").Value = 1102.5
$ws.Range("L99").Value = 4660
$ws.Range("M99").Value = 395.5
$ws.Range("N99").Value = -7656
$ws.Range("H105").Value = 3020
$ws.Range("I105").Value = 2087.647
$ws.Range("J105").Value = 5284.2856
$ws.Range("K105").Value = 2087.647
$ws.Range("L105").Value = 5284.2856
$ws.Range("M105").Value = -340.6469999999999
$ws.Range("N105").Value = -8778.285599999999
$ws.Range("H112").Value = 41200
$ws.Range("J112").Value = 41200
$ws.Range("L112").Value = 41200
$ws.Range("N112").Value = -44154

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2930.6
$ws.Range("I58").Value = 3100.75
$ws.Range("J58").Value = 2250
$ws.Range("K58").Value = 3100.75
$ws.Range("L58").Value = 2250
$ws.Range("M58").Value = -2897.75
$ws.Range("N58").Value = -2656
$ws.Range("H105").Value = 5489.1816
$ws.Range("I105").Value = 4862.625
$ws.Range("J105").Value = 7160
$ws.Range("K105").Value = 4862.625
$ws.Range("L105").Value = 7160
$ws.Range("M105").Value = -3115.625
$ws.Range("N105").Value = -10654
$ws.Range("H136").Value = 2930.6
$ws.Range("I136").Value = 3100.75
$ws.Range("J136").Value = 2250
$ws.Range("K136").Value = 9302.25
$ws.Range("L136").Value = 6750
$ws.Range("M136").Value = -6752.25
$ws.Range("N136").Value = -11850

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 121.333336
$ws.Range("I4").Value = 99
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 297
$ws.Range("L4").Value = 900
$ws.Range("M4").Value = -185
$ws.Range("N4").Value = -1124
$ws.Range("H121").Value = 703085.25
$ws.Range("I121").Value = 518.4
$ws.Range("J121").Value = 954002
$ws.Range("K121").Value = 1555.2
$ws.Range("L121").Value = 2862006
$ws.Range("M121").Value = -245.1999999999998
$ws.Range("N121").Value = -2864626
$ws.Range("H129").Value = 1489.9
$ws.Range("I129").Value = 754.9
$ws.Range("J129").Value = 1857.4
$ws.Range("K129").Value = 2264.7
$ws.Range("L129").Value = 5572.200000000001
$ws.Range("M129").Value = 2735.3
$ws.Range("N129").Value = -15572.2
$ws.Range("H132").Value = 695172
$ws.Range("I132").Value = 1316692.1
$ws.Range("J132").Value = 4594.1113
$ws.Range("K132").Value = 11850228.9
$ws.Range("L132").Value = 41347.00169999999
$ws.Range("M132").Value = -11847698.9
$ws.Range("N132").Value = -46407.00169999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2641.7
$ws.Range("I80").Value = 2440.6428
$ws.Range("J80").Value = 2817.625
$ws.Range("K80").Value = 2440.6428
$ws.Range("L80").Value = 2817.625
$ws.Range("M80").Value = -1442.6428
$ws.Range("N80").Value = -4813.625
$ws.Range("H83").Value = 2641.7
$ws.Range("I83").Value = 2440.6428
$ws.Range("J83").Value = 2817.625
$ws.Range("K83").Value = 12203.214
$ws.Range("L83").Value = 14088.125
$ws.Range("M83").Value = -7211.214
$ws.Range("N83").Value = -24072.125
$ws.Range("H102").Value = 2140.7083
$ws.Range("I102").Value = 1479.2354
$ws.Range("K102").Value = 1479.2354
$ws.Range("M102").Value = 142.7646

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2449.578
$ws.Range("I132").Value = 2327.634
$ws.Range("J132").Value = 3699.5
$ws.Range("K132").Value = 6982.902
$ws.Range("L132").Value = 11098.5
$ws.Range("M132").Value = -4452.902
$ws.Range("N132").Value = -16158.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6239.3335
$ws.Range("I62").Value = 3429
$ws.Range("K62").Value = 3429
$ws.Range("M62").Value = -2805
$ws.Range("H65").Value = 6239.3335
$ws.Range("I65").Value = 3429
$ws.Range("K65").Value = 17145
$ws.Range("M65").Value = -14025
$ws.Range("H122").Value = 8961.138000000001
$ws.Range("I122").Value = 10866.046
$ws.Range("J122").Value = 2974.2856
$ws.Range("K122").Value = 32598.138
$ws.Range("L122").Value = 8922.856800000001
$ws.Range("M122").Value = -30148.138
$ws.Range("N122").Value = -13822.8568
